# 2024-01 price update for 账户明细.xlsx
# Updates the "股价" (C column) on both sheets: 个人持仓 (sheet 1) and 家庭持仓 (sheet 2).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("个人持仓")
$ws2 = $wb.Worksheets.Item("家庭持仓")

# --- 个人持仓 ---
$ws1.Range("C2").Value = 42.78    # 中国平安
$ws1.Range("C3").Value = 29.17    # 达仁堂
$ws1.Range("C4").Value = 57.8     # 东阿阿胶
$ws1.Range("C5").Value = 49.42    # 云南白药
$ws1.Range("C6").Value = 29.08    # 白云山
$ws1.Range("C7").Value = 31.68    # 华东医药
$ws1.Range("C8").Value = 26.22    # 山东药玻
$ws1.Range("C9").Value = 24.51    # 三诺生物
$ws1.Range("C10").Value = 29.5    # 天坛生物
$ws1.Range("C11").Value = 98.5    # 洋河股份
$ws1.Range("C12").Value = 133.59  # 五粮液
$ws1.Range("C13").Value = 153.66  # 泸州老窖
$ws1.Range("C14").Value = 0.641   # 酒ETF
$ws1.Range("C15").Value = 8.57    # 恒顺醋业
$ws1.Range("C16").Value = 28.18   # 伊利股份
$ws1.Range("C17").Value = 28.57   # 双汇发展
$ws1.Range("C18").Value = 14.09   # 涪陵榨菜
$ws1.Range("C19").Value = 33.92   # 安琪酵母
$ws1.Range("C20").Value = 36.14   # 格力电器
$ws1.Range("C21").Value = 23.29   # 老板电器
$ws1.Range("C22").Value = 84.45   # 中国中免
$ws1.Range("C23").Value = 3.362   # 300ETF
$ws1.Range("C24").Value = 4.662   # 黄金ETF

# --- 家庭持仓 ---
$ws2.Range("C2").Value = 42.78    # 中国平安
$ws2.Range("C3").Value = 29.17    # 达仁堂
$ws2.Range("C4").Value = 57.8     # 东阿阿胶
$ws2.Range("C5").Value = 49.42    # 云南白药
$ws2.Range("C6").Value = 29.08    # 白云山
$ws2.Range("C7").Value = 31.68    # 华东医药
$ws2.Range("C8").Value = 26.22    # 山东药玻
$ws2.Range("C9").Value = 24.51    # 三诺生物
$ws2.Range("C10").Value = 29.5    # 天坛生物
$ws2.Range("C11").Value = 98.5    # 洋河股份
$ws2.Range("C12").Value = 153.66  # 泸州老窖
$ws2.Range("C13").Value = 133.59  # 五粮液
$ws2.Range("C14").Value = 225.2   # 山西汾酒
$ws2.Range("C15").Value = 0.641   # 酒ETF
$ws2.Range("C16").Value = 37.36   # 海天味业
$ws2.Range("C17").Value = 8.57    # 恒顺醋业
$ws2.Range("C18").Value = 28.18   # 伊利股份
$ws2.Range("C19").Value = 28.57   # 双汇发展
$ws2.Range("C20").Value = 14.09   # 涪陵榨菜
$ws2.Range("C21").Value = 33.92   # 安琪酵母
$ws2.Range("C22").Value = 36.14   # 格力电器
$ws2.Range("C23").Value = 23.29   # 老板电器
$ws2.Range("C24").Value = 84.45   # 中国中免
$ws2.Range("C25").Value = 3.362   # 300ETF
$ws2.Range("C26").Value = 0.853   # 100ETF
$ws2.Range("C27").Value = 4.662   # 黄金ETF
